# Mark the "PROMPT 4: isolation Matcher" prompt block as completed by
# applying strikethrough formatting to its body paragraphs (the heading
# itself and the following blank separator paragraph are left untouched),
# matching the pattern already used for the other completed prompts in
# this document.

$d = $word.ActiveDocument
$total = $d.Paragraphs.Count

# Locate the "PROMPT 4: isolation Matcher" Heading 3 paragraph.
$headingIndex = -1
for ($i = 1; $i -le $total; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 3" -and $p.Range.Text -like "PROMPT 4:*isolation*") {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not find the 'PROMPT 4: isolation Matcher' heading paragraph."
}

# Find the next Heading 3 (start of the following prompt) so we know
# where this prompt's block ends.
$nextHeadingIndex = $total + 1
for ($i = $headingIndex + 1; $i -le $total; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Heading 3") {
        $nextHeadingIndex = $i
        break
    }
}

# The block's body paragraphs run from right after the heading up to
# (but excluding) the trailing blank "Normal" separator paragraph that
# precedes the next heading.
$firstBodyIndex = $headingIndex + 1
$lastBodyIndex = $nextHeadingIndex - 1
while ($lastBodyIndex -ge $firstBodyIndex -and $d.Paragraphs.Item($lastBodyIndex).Style.NameLocal -eq "Normal") {
    $lastBodyIndex = $lastBodyIndex - 1
}

for ($i = $firstBodyIndex; $i -le $lastBodyIndex; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "Body Text") {
        # Apply to the whole paragraph range, including its trailing
        # paragraph mark, so the mark's own run properties also pick up
        # strikethrough (mirrors the other completed prompts).
        $p.Range.Font.StrikeThrough = 1
    } else {
        # For paragraphs whose mark formatting should stay untouched
        # (e.g. the "First Paragraph" styled intro line), only strike the
        # visible text, excluding the trailing paragraph-mark character.
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.Font.StrikeThrough = 1
    }
}

Write-Output "Struck through paragraphs $firstBodyIndex..$lastBodyIndex (heading at $headingIndex)"
